$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1514.1428
$ws.Range("I28").Value = 1514.1428
$ws.Range("K28").Value = 1514.1428
$ws.Range("M28").Value = -1029.1428
$ws.Range("H70").Value = 2509.56
$ws.Range("I70").Value = 1991.0526
$ws.Range("J70").Value = 4151.5
$ws.Range("K70").Value = 5973.1578
$ws.Range("L70").Value = 12454.5
$ws.Range("M70").Value = -5703.1578
$ws.Range("N70").Value = -12994.5
$ws.Range("H73").Value = 2509.56
$ws.Range("I73").Value = 1991.0526
$ws.Range("J73").Value = 4151.5
$ws.Range("K73").Value = 5973.1578
$ws.Range("L73").Value = 12454.5
$ws.Range("M73").Value = -5037.1578
$ws.Range("N73").Value = -14326.5
$ws.Range("H137").Value = 605445
$ws.Range("I137").Value = 1363229.1
$ws.Range("J137").Value = 2662.0908
$ws.Range("K137").Value = 4089687.3
$ws.Range("L137").Value = 7986.2724
$ws.Range("M137").Value = -4087137.3
$ws.Range("N137").Value = -13086.2724
$ws.Range("H140").Value = 72705.766
$ws.Range("J140").Value = 72705.766
$ws.Range("L140").Value = 72705.766
$ws.Range("N140").Value = -83065.766
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4074.1594
$ws.Range("I32").Value = 3759.2456
$ws.Range("K32").Value = 3759.2456
$ws.Range("M32").Value = -3472.2456
$ws.Range("H74").Value = 248975.38
$ws.Range("I74").Value = 395310.8
$ws.Range("J74").Value = 2034.3125
$ws.Range("K74").Value = 395310.8
$ws.Range("L74").Value = 2034.3125
$ws.Range("M74").Value = -394436.8
$ws.Range("N74").Value = -3782.3125
$ws.Range("H77").Value = 248975.38
$ws.Range("I77").Value = 395310.8
$ws.Range("J77").Value = 2034.3125
$ws.Range("K77").Value = 1976554
$ws.Range("L77").Value = 10171.5625
$ws.Range("M77").Value = -1972186
$ws.Range("N77").Value = -18907.5625
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents() | Out-Null
$ws.Range("H132").Value = 2501.697
$ws.Range("I132").Value = 1551.8334
$ws.Range("J132").Value = 5034.6665
$ws.Range("K132").Value = 4655.5002
$ws.Range("L132").Value = 15103.9995
$ws.Range("M132").Value = -2125.5002
$ws.Range("N132").Value = -20163.9995
$ws.Range("H137").Value = 40336
$ws.Range("J137").Value = 40336
$ws.Range("L137").Value = 40336
$ws.Range("N137").Value = -50536
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4200.4
$ws.Range("J4").Value = 4200.4
$ws.Range("L4").Value = 4200.4
$ws.Range("N4").Value = -4424.4
$ws.Range("H31").Value = 243591.1
$ws.Range("I31").Value = 541321.2
$ws.Range("J31").Value = 3486.1936
$ws.Range("K31").Value = 541321.2
$ws.Range("L31").Value = 3486.1936
$ws.Range("M31").Value = -541026.2
$ws.Range("N31").Value = -4076.1936
$ws.Range("H34").Value = 243591.1
$ws.Range("I34").Value = 541321.2
$ws.Range("J34").Value = 3486.1936
$ws.Range("K34").Value = 541321.2
$ws.Range("L34").Value = 3486.1936
$ws.Range("M34").Value = -541119.2
$ws.Range("N34").Value = -3890.1936
$ws.Range("H58").Value = 2782.2188
$ws.Range("I58").Value = 1501.24
$ws.Range("J58").Value = 7357.143
$ws.Range("K58").Value = 1501.24
$ws.Range("L58").Value = 7357.143
$ws.Range("M58").Value = -1298.24
$ws.Range("N58").Value = -7763.143
$ws.Range("H62").Value = 4333.3335
$ws.Range("I62").Value = 4333.3335
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4333.3335
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3709.3335
$ws.Range("N62").ClearContents() | Out-Null
$ws.Range("H65").Value = 4333.3335
$ws.Range("I65").Value = 4333.3335
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21666.6675
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18546.6675
$ws.Range("N65").ClearContents() | Out-Null
$ws.Range("H136").Value = 2782.2188
$ws.Range("I136").Value = 1501.24
$ws.Range("J136").Value = 7357.143
$ws.Range("K136").Value = 4503.72
$ws.Range("L136").Value = 22071.429
$ws.Range("M136").Value = -1953.72
$ws.Range("N136").Value = -27171.429
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1608208.2
$ws.Range("I4").Value = 8036541.5
$ws.Range("J4").Value = 1125
$ws.Range("K4").Value = 24109624.5
$ws.Range("L4").Value = 3375
$ws.Range("M4").Value = -24109512.5
$ws.Range("N4").Value = -3599
$ws.Range("H107").Value = 12816.19
$ws.Range("I107").Value = 389.23077
$ws.Range("J107").Value = 23586.223
$ws.Range("K107").Value = 1167.69231
$ws.Range("L107").Value = 70758.66900000001
$ws.Range("M107").Value = 752.3076900000001
$ws.Range("N107").Value = -74598.66900000001
$ws.Range("H131").Value = 6494435.5
$ws.Range("I131").Value = 83334150
$ws.Range("J131").Value = 938.39435
$ws.Range("K131").Value = 250002450
$ws.Range("L131").Value = 2815.18305
$ws.Range("M131").Value = -249997410
$ws.Range("N131").Value = -12895.18305
$ws.Range("H132").Value = 2789.077
$ws.Range("I132").Value = 950
$ws.Range("J132").Value = 3606.4443
$ws.Range("K132").Value = 8550
$ws.Range("L132").Value = 32457.9987
$ws.Range("M132").Value = -6020
$ws.Range("N132").Value = -37517.9987
$ws.Range("H133").Value = 6704.2856
$ws.Range("I133").Value = 7021.6665
$ws.Range("J133").Value = 4800
$ws.Range("K133").Value = 21064.9995
$ws.Range("L133").Value = 14400
$ws.Range("M133").Value = -16004.9995
$ws.Range("N133").Value = -24520
$ws.Range("H134").Value = 1618
$ws.Range("I134").Value = 1618
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4854
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 216
$ws.Range("N134").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2220.566
$ws.Range("I132").Value = 1160.5186
$ws.Range("K132").Value = 3481.5558
$ws.Range("M132").Value = -951.5558000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 951.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents() | Out-Null
